$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Benzin" row (row 6) entirely - shifts rows 7:12 up to 6:11
$ws.Rows("6:6").Delete()

# Add two new commodity rows at the bottom of the table ("Holz" and "Rind"),
# same category/indicator/value pattern as the preceding "Kakao" row (now row 11)
$ws.Range("A12").Value = "Holz"
$ws.Range("B12").Value = "Energie u. weiteres"
$ws.Range("D12").Value = "Rohstoff Indikator"
$ws.Range("E12").Value = 3053

$ws.Range("A13").Value = "Rind"
$ws.Range("B13").Value = "Energie u. weiteres"
$ws.Range("D13").Value = "Rohstoff Indikator"
$ws.Range("E13").Value = 3053

# Update the active selection to reflect the new last row
$ws.Range("A13").Select()
